{"js": "// Bump the document version string from \"Version 1.\" to \"Version 2.\"\n// (the \"wireframes version 2\" edit referenced by the commit message).\nconst body = context.document.body;\n\n// Find the \"1.\" that follows \"Version \" and swap it for \"2.\".\nconst results = body.search(\"1.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"2.\", Word.InsertLocation.replace);\n} else {\n  // Fallback: if the exact \"1.\" substring isn't found (e.g. already edited),\n  // replace the whole paragraph text directly.\n  const results2 = body.search(\"Version 1.\", { matchCase: true });\n  results2.load(\"items\");\n  await context.sync();\n  if (results2.items.length > 0) {\n    results2.items[0].insertText(\"Version 2.\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Bump the document version string from \"Version 1.\" to \"Version 2.\"\n# (the \"wireframes version 2\" edit referenced by the commit message).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"1.\", $false, $false, $false, $false, $false, $true, 1, $false, \"2.\", 2)\n\nif (-not $found) {\n    # Fallback: exact \"1.\" substring not found (e.g. already edited) -\n    # replace the whole paragraph text directly.\n    $range2 = $d.Content\n    $range2.Find.Execute(\"Version 1.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version 2.\", 2)\n}\n"}
